# Update "want to go" counts (column F) across the four sheets to match
# the newly scraped numbers from the generator run (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1300
$ws1.Range("F7").Value = 386
$ws1.Range("F10").Value = 10259
$ws1.Range("F11").Value = 89
$ws1.Range("F18").Value = 287
$ws1.Range("F26").Value = 56
$ws1.Range("F30").Value = 53
$ws1.Range("F31").Value = 570
$ws1.Range("F40").Value = 15
$ws1.Range("F41").Value = 124
$ws1.Range("F43").Value = 332
$ws1.Range("F44").Value = 75
$ws1.Range("F47").Value = 28
$ws1.Range("F49").Value = 63

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 64

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 336

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 336
$ws4.Range("F9").Value = 1300
$ws4.Range("F15").Value = 10259
$ws4.Range("F19").Value = 287
$ws4.Range("F23").Value = 56
$ws4.Range("F25").Value = 53
$ws4.Range("F27").Value = 570
$ws4.Range("F32").Value = 64
$ws4.Range("F37").Value = 124
$ws4.Range("F40").Value = 332
$ws4.Range("F41").Value = 75
$ws4.Range("F47").Value = 28
$ws4.Range("F49").Value = 63
